$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-23 (name, total_registros), reflecting the updated
# monitoring snapshot. Several rows were re-sorted because the sheet is
# sorted descending by total_registros.
$data = @(
    @("FARFAN MONTOYA ROSSANA ISABEL", 78),
    @("DAVILA CORDOVA MARIBEL", 69),
    @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 64),
    @("SANCHEZ ULLOA CESAR AUGUSTO", 63),
    @("RAMOS RAMOS HANDY JAIR", 62),
    @("GARCIA GUTIERREZ LUIS ARTURO", 61),
    @("FIORELA KEILY GUTIERREZ CRUZ", 59),
    @("CARBAJAL RAMOS JESUS MARINA", 58),
    @("GONZALES VALLE SEBASTIAN", 58),
    @("DELGADO DELGADO RONI", 57),
    @("VERDE LIZARRAGA DEYSI EUFEMIA", 56),
    @("OLIVA ALVA GOSSELYN NASSIRA", 55),
    @("JOSSY IVANA SUÁREZ ZAVALETA", 55),
    @("BAZAN TEJADA JOSE VICENTE", 54),
    @("DANY DARWIN VILLACORTA SAAVEDRA", 53),
    @("CASTILLO QUEZADA DIEGO ALONSO", 50),
    @("ARANEDA LOPEZ MARCO VIERI", 50),
    @("GUERRA CALDERON ESTHEFANY NICOLLE", 48),
    @("CARDENAS CAMPOJO MARY PAULA", 46),
    @("JAVE CHAVEZ ANGHELO MARTIN", 45),
    @("ALVITES CAMPOS SERGIO MARTIN", 22),
    @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
